$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '67.531.21'
$ws.Range('E2').Value = '  -1.95%  '
Set-TextValue 'D3' '2.425.06'
$ws.Range('E3').Value = '  -1.64%  '
$ws.Range('E4').Value = '  +0.08%  '
Set-TextValue 'D5' '550.89'
$ws.Range('E5').Value = '  -1.29%  '
Set-TextValue 'D6' '158.81'
$ws.Range('E6').Value = '  -2.10%  '
$ws.Range('E7').Value = '  +0.10%  '
Set-TextValue 'D8' '0.506'
$ws.Range('E8').Value = '  +0.42%  '
$ws.Range('E9').Value = '  +5.41%  '
$ws.Range('E10').Value = '  -0.96%  '
$ws.Range('E11').Value = '  -1.56%  '
$ws.Range('E12').Value = '  -0.82%  '
Set-TextValue 'D13' '67.693.51'
$ws.Range('E13').Value = '  -1.47%  '
Set-TextValue 'D14' '0.0000169'
$ws.Range('E14').Value = '  +0.27%  '
Set-TextValue 'D15' '22.92'
$ws.Range('E15').Value = '  -2.42%  '
Set-TextValue 'D16' '10.32'
$ws.Range('E16').Value = '  -3.89%  '
Set-TextValue 'D17' '328.58'
$ws.Range('E17').Value = '  -3.86%  '
Set-TextValue 'D18' '6.83'
$ws.Range('E18').Value = '  -2.93%  '
Set-TextValue 'D19' '3.78'
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('E20').Value = '  -0.39%  '
Set-TextValue 'D21' '1.83'
$ws.Range('E21').Value = '  -1.38%  '
Set-TextValue 'D22' '65.97'
$ws.Range('E22').Value = '  -1.23%  '
Set-TextValue 'D23' '3.60'
$ws.Range('E23').Value = '  -1.88%  '
Set-TextValue 'D24' '8.03'
$ws.Range('E24').Value = '  -1.15%  '
Set-TextValue 'D25' '0.0₃0798'
$ws.Range('E25').Value = '  -2.43%  '
$ws.Range('E26').Value = '  -2.23%  '
Set-TextValue 'D27' '1.00'
$ws.Range('E27').Value = '  +0.07%  '
Set-TextValue 'D28' '413.93'
$ws.Range('E28').Value = '  -5.34%  '
$ws.Range('E29').Value = '  -1.79%  '
$ws.Range('E30').Value = '  -1.44%  '
Set-TextValue 'D31' '159.31'
$ws.Range('E31').Value = '  +1.32%  '
$ws.Range('E32').Value = '  -0.60%  '
$ws.Range('E33').Value = '  -0.12%  '
$ws.Range('E34').Value = '  -0.46%  '
Set-TextValue 'D35' '0.104'
$ws.Range('E35').Value = '  -4.49%  '
Set-TextValue 'D36' '0.293'
$ws.Range('E36').Value = '  -3.09%  '
Set-TextValue 'D37' '4.23'
$ws.Range('E37').Value = '  -4.70%  '
$ws.Range('E38').Value = '  -1.80%  '
Set-TextValue 'D39' '1.06'
$ws.Range('E39').Value = '  -3.11%  '
Set-TextValue 'D40' '1.98'
$ws.Range('E40').Value = '  -4.20%  '
Set-TextValue 'D41' '3.30'
$ws.Range('E41').Value = '  -1.24%  '
Set-TextValue 'D42' '129.66'
$ws.Range('E42').Value = '  -2.25%  '
$ws.Range('E43').Value = '  -1.14%  '
Set-TextValue 'D44' '0.475'
$ws.Range('E44').Value = '  -1.72%  '
Set-TextValue 'D45' '0.552'
$ws.Range('E45').Value = '  -1.81%  '
Set-TextValue 'D46' '0.0912'
$ws.Range('E46').Value = '  +0.58%  '
$ws.Range('E47').Value = '  +0.33%  '
$ws.Range('E48').Value = '  -8.08%  '
Set-TextValue 'D49' '16.43'
$ws.Range('E49').Value = '  -2.51%  '
Set-TextValue 'D50' '0.0₆0202'
$ws.Range('E50').Value = '  +0.79%  '
Set-TextValue 'D51' '0.0427'
$ws.Range('E51').Value = '  -0.60%  '
